# Enhance 3D Gun Barrel Plot functionality by adding New Mexico land survey
# system support, increasing plot size, and updating workflow tasks.
#
# This script updates the stratigraphic color lookup table: the color swatch
# used for the "Wolfcamp XY" row is changed to a light gray (#D3D3D3, a new
# color used to flag the New Mexico land survey system), and the swatch used
# for "Wolfcamp D" is changed to a light blue (#ADD8E6). The previous colors
# for those two rows (#4D8FD1 and #5EB85E) are no longer referenced anywhere
# in the sheet, so they drop out of the shared string table, while the two
# new colors are appended to it.
#
# It also updates the current selection on the sheet (from H16 to E33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wolfcamp XY row (row 14) - new "flag" color for NM land survey system.
$ws.Range("H14").Value2 = "#D3D3D3"

# Wolfcamp D row (row 20) - new plot color.
$ws.Range("H20").Value2 = "#ADD8E6"

# Update the active selection/cell.
$ws.Range("E33").Select()
